# Fix result structure and parameter extraction
#
# Column B ("Parameter") previously held a long, verbose sentence (often
# duplicated from column D's "Details" text). This sets column B to just
# the short parameter name, which is already present (single-quoted) at
# the start of each row's Details (column D) text, e.g.
#   D2 = "Fail: 's:event:type' parameter missing. URL: ... Please add the 's:event:type' parameter."
#   -> B2 = "s:event:type"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $details = $ws.Cells.Item($row, 4).Value2
    if ($details -eq $null) { continue }

    $details = [string]$details

    $firstQuote = $details.IndexOf("'")
    if ($firstQuote -ge 0) {
        $secondQuote = $details.IndexOf("'", $firstQuote + 1)
        if ($secondQuote -gt $firstQuote) {
            $paramName = $details.Substring($firstQuote + 1, $secondQuote - $firstQuote - 1)
            $ws.Cells.Item($row, 2).Value = $paramName
        }
    }
}
